# Apply scheduled market-data updates to the leve profit tables across sheets
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 10
$ws_ALC.Range("H10").Value = 15238
$ws_ALC.Range("I10").Value = 18747.5
$ws_ALC.Range("J10").Value = 1200
$ws_ALC.Range("K10").Value = 18747.5
$ws_ALC.Range("L10").Value = 1200
$ws_ALC.Range("M10").Value = -18454.5
$ws_ALC.Range("N10").Value = -1786

# ALC row 64
$ws_ALC.Range("H64").Value = 3520.75
$ws_ALC.Range("I64").Value = 3571.1428
$ws_ALC.Range("J64").Value = 3500
$ws_ALC.Range("K64").Value = 3571.1428
$ws_ALC.Range("L64").Value = 3500
$ws_ALC.Range("M64").Value = -3323.1428
$ws_ALC.Range("N64").Value = -3996

# ALC row 67
$ws_ALC.Range("H67").Value = 3520.75
$ws_ALC.Range("I67").Value = 3571.1428
$ws_ALC.Range("J67").Value = 3500
$ws_ALC.Range("K67").Value = 3571.1428
$ws_ALC.Range("L67").Value = 3500
$ws_ALC.Range("M67").Value = -2713.1428
$ws_ALC.Range("N67").Value = -5216

# ALC row 86
$ws_ALC.Range("H86").Value = 3000
$ws_ALC.Range("I86").Value = 3000
$ws_ALC.Range("J86").Value = 0
$ws_ALC.Range("K86").Value = 3000
$ws_ALC.Range("L86").Value = 0
$ws_ALC.Range("M86").Value = -1877
$ws_ALC.Range("N86").ClearContents()

# ALC row 89
$ws_ALC.Range("H89").Value = 3000
$ws_ALC.Range("I89").Value = 3000
$ws_ALC.Range("J89").Value = 0
$ws_ALC.Range("K89").Value = 15000
$ws_ALC.Range("L89").Value = 0
$ws_ALC.Range("M89").Value = -9384
$ws_ALC.Range("N89").ClearContents()

# ALC row 111
$ws_ALC.Range("H111").Value = 1166
$ws_ALC.Range("I111").Value = 0
$ws_ALC.Range("J111").Value = 1166
$ws_ALC.Range("K111").Value = 0
$ws_ALC.Range("L111").Value = 3498
$ws_ALC.Range("N111").Value = -9632
$ws_ALC.Range("M111").ClearContents()

# ARM row 8
$ws_ARM.Range("H8").Value = 0
$ws_ARM.Range("J8").Value = 0
$ws_ARM.Range("L8").Value = 0
$ws_ARM.Range("N8").ClearContents()

# ARM row 74
$ws_ARM.Range("H74").Value = 1422.2222
$ws_ARM.Range("I74").Value = 825
$ws_ARM.Range("K74").Value = 825
$ws_ARM.Range("M74").Value = 49

# ARM row 77
$ws_ARM.Range("H77").Value = 1422.2222
$ws_ARM.Range("I77").Value = 825
$ws_ARM.Range("K77").Value = 4125
$ws_ARM.Range("M77").Value = 243

# CRP row 4
$ws_CRP.Range("H4").Value = 10000
$ws_CRP.Range("I4").Value = 0
$ws_CRP.Range("K4").Value = 0
$ws_CRP.Range("M4").ClearContents()

# CRP row 7
$ws_CRP.Range("H7").Value = 138.55556
$ws_CRP.Range("I7").Value = 90.90000000000001
$ws_CRP.Range("J7").Value = 198.125
$ws_CRP.Range("K7").Value = 90.90000000000001
$ws_CRP.Range("L7").Value = 198.125
$ws_CRP.Range("M7").Value = 22.09999999999999
$ws_CRP.Range("N7").Value = -424.125

# CRP row 16
$ws_CRP.Range("H16").Value = 1641.4286
$ws_CRP.Range("I16").Value = 1537.8
$ws_CRP.Range("J16").Value = 1900.5
$ws_CRP.Range("K16").Value = 1537.8
$ws_CRP.Range("L16").Value = 1900.5
$ws_CRP.Range("M16").Value = -1250.8
$ws_CRP.Range("N16").Value = -2474.5

# CRP row 62
$ws_CRP.Range("H62").Value = 52730
$ws_CRP.Range("I62").Value = 73785.71000000001
$ws_CRP.Range("J62").Value = 3600
$ws_CRP.Range("K62").Value = 73785.71000000001
$ws_CRP.Range("L62").Value = 3600
$ws_CRP.Range("M62").Value = -73161.71000000001
$ws_CRP.Range("N62").Value = -4848

# CRP row 65
$ws_CRP.Range("H65").Value = 52730
$ws_CRP.Range("I65").Value = 73785.71000000001
$ws_CRP.Range("J65").Value = 3600
$ws_CRP.Range("K65").Value = 368928.55
$ws_CRP.Range("L65").Value = 18000
$ws_CRP.Range("M65").Value = -365808.55
$ws_CRP.Range("N65").Value = -24240

# CRP row 99
$ws_CRP.Range("H99").Value = 2279.3
$ws_CRP.Range("I99").Value = 2356.8
$ws_CRP.Range("J99").Value = 2201.8
$ws_CRP.Range("K99").Value = 2356.8
$ws_CRP.Range("L99").Value = 2201.8
$ws_CRP.Range("M99").Value = -858.8000000000002
$ws_CRP.Range("N99").Value = -5197.8

# CRP row 113
$ws_CRP.Range("H113").Value = 1641.4286
$ws_CRP.Range("I113").Value = 1537.8
$ws_CRP.Range("J113").Value = 1900.5
$ws_CRP.Range("K113").Value = 1537.8
$ws_CRP.Range("L113").Value = 1900.5
$ws_CRP.Range("M113").Value = 632.2
$ws_CRP.Range("N113").Value = -6240.5

# CRP row 126
$ws_CRP.Range("H126").Value = 2279.3
$ws_CRP.Range("I126").Value = 2356.8
$ws_CRP.Range("J126").Value = 2201.8
$ws_CRP.Range("K126").Value = 7070.400000000001
$ws_CRP.Range("L126").Value = 6605.400000000001
$ws_CRP.Range("M126").Value = -4600.400000000001
$ws_CRP.Range("N126").Value = -11545.4

# CUL row 16
$ws_CUL.Range("H16").Value = 1500
$ws_CUL.Range("I16").Value = 750
$ws_CUL.Range("J16").Value = 3000
$ws_CUL.Range("K16").Value = 2250
$ws_CUL.Range("L16").Value = 9000
$ws_CUL.Range("M16").Value = -2077
$ws_CUL.Range("N16").Value = -9346

# CUL row 23
$ws_CUL.Range("H23").Value = 344.5
$ws_CUL.Range("I23").Value = 212
$ws_CUL.Range("J23").Value = 373.30435
$ws_CUL.Range("K23").Value = 636
$ws_CUL.Range("L23").Value = 1119.91305
$ws_CUL.Range("M23").Value = -401
$ws_CUL.Range("N23").Value = -1589.91305

# CUL row 117
$ws_CUL.Range("H117").Value = 759.5
$ws_CUL.Range("J117").Value = 780.9231
$ws_CUL.Range("L117").Value = 2342.7693
$ws_CUL.Range("N117").Value = -9226.7693

# CUL row 131
$ws_CUL.Range("H131").Value = 10010297
$ws_CUL.Range("I131").Value = 432.5
$ws_CUL.Range("J131").Value = 10427375
$ws_CUL.Range("K131").Value = 1297.5
$ws_CUL.Range("L131").Value = 31282125
$ws_CUL.Range("M131").Value = 3742.5
$ws_CUL.Range("N131").Value = -31292205

# GSM row 10
$ws_GSM.Range("H10").Value = 13960
$ws_GSM.Range("J10").Value = 13960
$ws_GSM.Range("L10").Value = 13960
$ws_GSM.Range("N10").Value = -14298

# GSM row 102
$ws_GSM.Range("H102").Value = 2112.6333
$ws_GSM.Range("I102").Value = 1747.0454
$ws_GSM.Range("J102").Value = 3118
$ws_GSM.Range("K102").Value = 1747.0454
$ws_GSM.Range("L102").Value = 3118
$ws_GSM.Range("M102").Value = -125.0454
$ws_GSM.Range("N102").Value = -6362

# GSM row 107
$ws_GSM.Range("H107").Value = 3792.6667
$ws_GSM.Range("I107").Value = 4391.2
$ws_GSM.Range("J107").Value = 800
$ws_GSM.Range("K107").Value = 4391.2
$ws_GSM.Range("L107").Value = 800
$ws_GSM.Range("M107").Value = -2471.2
$ws_GSM.Range("N107").Value = -4640

# GSM row 113
$ws_GSM.Range("H113").Value = 1567.5714
$ws_GSM.Range("I113").Value = 925.3333
$ws_GSM.Range("K113").Value = 925.3333
$ws_GSM.Range("M113").Value = 1244.6667

# GSM row 118
$ws_GSM.Range("H118").Value = 0
$ws_GSM.Range("J118").Value = 0
$ws_GSM.Range("L118").Value = 0
$ws_GSM.Range("N118").ClearContents()

# GSM row 120
$ws_GSM.Range("H120").Value = 30000
$ws_GSM.Range("J120").Value = 30000
$ws_GSM.Range("L120").Value = 30000
$ws_GSM.Range("N120").Value = -39676

# GSM row 121
$ws_GSM.Range("H121").Value = 22900
$ws_GSM.Range("J121").Value = 22900
$ws_GSM.Range("L121").Value = 22900
$ws_GSM.Range("N121").Value = -26394

# LTW row 16
$ws_LTW.Range("H16").Value = 1062.88
$ws_LTW.Range("I16").Value = 868.1739
$ws_LTW.Range("K16").Value = 868.1739
$ws_LTW.Range("M16").Value = -698.1739

# LTW row 40
$ws_LTW.Range("H40").Value = 5529.1665
$ws_LTW.Range("I40").Value = 6070
$ws_LTW.Range("J40").Value = 4447.5
$ws_LTW.Range("K40").Value = 6070
$ws_LTW.Range("L40").Value = 4447.5
$ws_LTW.Range("M40").Value = -5934
$ws_LTW.Range("N40").Value = -4719.5

# LTW row 100
$ws_LTW.Range("H100").Value = 7622.222
$ws_LTW.Range("I100").Value = 17566.666
$ws_LTW.Range("J100").Value = 2650
$ws_LTW.Range("K100").Value = 17566.666
$ws_LTW.Range("L100").Value = 2650
$ws_LTW.Range("M100").Value = -17025.666
$ws_LTW.Range("N100").Value = -3732

# LTW row 122
$ws_LTW.Range("H122").Value = 22506150
$ws_LTW.Range("I122").Value = 31255312
$ws_LTW.Range("J122").Value = 16673375
$ws_LTW.Range("K122").Value = 93765936
$ws_LTW.Range("L122").Value = 50020125
$ws_LTW.Range("M122").Value = -93763486
$ws_LTW.Range("N122").Value = -50025025

# WVR row 100
$ws_WVR.Range("H100").Value = 20453
$ws_WVR.Range("I100").Value = 28895.715
$ws_WVR.Range("J100").Value = 753.3333
$ws_WVR.Range("K100").Value = 57791.43
$ws_WVR.Range("L100").Value = 1506.6666
$ws_WVR.Range("M100").Value = -57250.43
$ws_WVR.Range("N100").Value = -2588.6666
